$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'36.933.72"
$ws.Range("E2").Formula = "'  -0.29%  "
$ws.Range("D3").Formula = "'2.041.51"
$ws.Range("E3").Formula = "'  -0.52%  "
$ws.Range("E4").Formula = "'  +0.06%  "
$ws.Range("D5").Formula = "'245.94"
$ws.Range("E5").Formula = "'  -1.20%  "
$ws.Range("D6").Formula = "'0.657"
$ws.Range("E6").Formula = "'  -0.95%  "
$ws.Range("D7").Formula = "'57.60"
$ws.Range("E7").Formula = "'  -1.19%  "
$ws.Range("E8").Formula = "'  -0.05%  "
$ws.Range("D9").Formula = "'0.378"
$ws.Range("E9").Formula = "'  -1.14%  "
$ws.Range("D10").Formula = "'0.0767"
$ws.Range("E10").Formula = "'  -1.95%  "
$ws.Range("D11").Formula = "'0.110"
$ws.Range("E11").Formula = "'  +1.81%  "
$ws.Range("D12").Formula = "'15.54"
$ws.Range("E12").Formula = "'  -1.96%  "
$ws.Range("D13").Formula = "'0.888"
$ws.Range("E13").Formula = "'  +11.19%  "
$ws.Range("D14").Formula = "'2.340.88"
$ws.Range("E14").Formula = "'  -0.35%  "
$ws.Range("D15").Formula = "'5.67"
$ws.Range("E15").Formula = "'  +1.10%  "
$ws.Range("D16").Formula = "'2.071.38"
$ws.Range("E16").Formula = "'  +0.87%  "
$ws.Range("D17").Formula = "'18.14"
$ws.Range("E17").Formula = "'  +9.33%  "
$ws.Range("D18").Formula = "'36.940.83"
$ws.Range("E18").Formula = "'  -0.27%  "
$ws.Range("D19").Formula = "'73.85"
$ws.Range("E19").Formula = "'  -1.14%  "
$ws.Range("D20").Formula = "'0.0₃0892"
$ws.Range("E20").Formula = "'  -0.24%  "
$ws.Range("D21").Formula = "'5.39"
$ws.Range("E21").Formula = "'  +0.71%  "
$ws.Range("D22").Formula = "'235.57"
$ws.Range("E22").Formula = "'  -0.31%  "
$ws.Range("E23").Formula = "'  +0.03%  "
$ws.Range("E24").Formula = "'  +2.61%  "
$ws.Range("D25").Formula = "'170.07"
$ws.Range("E25").Formula = "'  +1.17%  "
$ws.Range("D26").Formula = "'9.45"
$ws.Range("E26").Formula = "'  +3.28%  "
$ws.Range("D27").Formula = "'2.13"
$ws.Range("E27").Formula = "'  -3.94%  "
$ws.Range("D28").Formula = "'19.87"
$ws.Range("E28").Formula = "'  +0.14%  "
$ws.Range("D29").Formula = "'5.45"
$ws.Range("E29").Formula = "'  +15.92%  "
$ws.Range("E30").Formula = "'  -0.92%  "
$ws.Range("E31").Formula = "'  -0.16%  "
$ws.Range("D32").Formula = "'4.74"
$ws.Range("E32").Formula = "'  +6.29%  "
$ws.Range("D33").Formula = "'0.0614"
$ws.Range("E33").Formula = "'  -0.27%  "
$ws.Range("E34").Formula = "'  +0.08%  "
$ws.Range("D35").Formula = "'0.0869"
$ws.Range("E35").Formula = "'  -2.03%  "
$ws.Range("E36").Formula = "'  +6.07%  "
$ws.Range("D37").Formula = "'2.24"
$ws.Range("E37").Formula = "'  +0.41%  "
$ws.Range("D38").Formula = "'1.31"
$ws.Range("E38").Formula = "'  -1.97%  "
$ws.Range("E39").Formula = "'  -1.77%  "
$ws.Range("D40").Formula = "'5.07"
$ws.Range("E40").Formula = "'  +2.97%  "
$ws.Range("D41").Formula = "'0.0984"
$ws.Range("E41").Formula = "'  -8.29%  "
$ws.Range("E42").Formula = "'  +0.69%  "
$ws.Range("E43").Formula = "'  +1.45%  "
$ws.Range("B44").Formula = "'InjectiveProtocol"
$ws.Range("C44").Formula = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Formula = "'17.00"
$ws.Range("E44").Formula = "'  -2.35%  "
$ws.Range("B45").Formula = "'Aave"
$ws.Range("C45").Formula = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Formula = "'96.79"
$ws.Range("E45").Formula = "'  +0.72%  "
$ws.Range("D46").Formula = "'1.291.52"
$ws.Range("E46").Formula = "'  +0.56%  "
$ws.Range("D47").Formula = "'2.35"
$ws.Range("E47").Formula = "'  -4.50%  "
$ws.Range("B48").Formula = "'FTXToken"
$ws.Range("C48").Formula = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Formula = "'3.77"
$ws.Range("E48").Formula = "'  +7.30%  "
$ws.Range("B49").Formula = "'MXToken"
$ws.Range("C49").Formula = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Formula = "'2.85"
$ws.Range("E49").Formula = "'  -0.03%  "
$ws.Range("D50").Formula = "'6.73"
$ws.Range("E50").Formula = "'  +0.60%  "
$ws.Range("D51").Formula = "'2.226.90"
$ws.Range("E51").Formula = "'  -0.56%  "
